$d = $word.ActiveDocument

$replacements = @(
    @("99×63=6237", "12×61=732"),
    @("69×53=3657", "48×56=2688"),
    @("85×76=6460", "43×58=2494"),
    @("14×84=1176", "92×39=3588"),
    @("81×64=5184", "63×61=3843"),
    @("42×42=1764", "23×44=1012"),
    @("54×95=5130", "82×14=1148"),
    @("63×40=2520", "23×91=2093"),
    @("85×37=3145", "63×95=5985"),
    @("56×54=3024", "27×26=702"),
    @("54×49=2646", "55×83=4565"),
    @("17×34=578", "93×13=1209"),
    @("73×55=4015", "71×59=4189"),
    @("30×50=1500", "37×39=1443"),
    @("60×92=5520", "59×12=708"),
    @("72×99=7128", "38×38=1444"),
    @("19×61=1159", "64×80=5120"),
    @("45×45=2025", "35×25=875"),
    @("70×20=1400", "16×23=368"),
    @("49×92=4508", "12×62=744"),
    @("87×79=6873", "26×36=936"),
    @("44×50=2200", "45×95=4275"),
    @("97×73=7081", "71×73=5183"),
    @("92×91=8372", "61×17=1037"),
    @("21×80=1680", "75×21=1575")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
